# The "Date" column (column D, rows 2-6) holds plain text values like
# "09/29/2017 00:00:00" instead of "09/29/2017" - strip the redundant
# " 00:00:00" midnight timestamp so only the date portion remains, while
# keeping the cell content as literal text (not a real Excel date/number).
#
# Because these are text cells, we must pre-format each cell as Text
# ("@") before writing back the trimmed string - otherwise Excel's usual
# Range.Value smart-parsing would recognize "MM/DD/YYYY" as a date and
# silently convert the cell into a numeric date serial value.
#
# Demonstrates two different ways of writing into the worksheet, per cell:
#   1) Range("A1").Value
#   2) Cells.Item(row, column).Value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Strip-MidnightTimestamp([string]$text) {
    return ($text -replace " 00:00:00$", "")
}

# --- Method 1: Range(...).Value -------------------------------------------
$range = $ws.Range("D2")
$range.NumberFormat = "@"
$range.Value = Strip-MidnightTimestamp $range.Text

$range = $ws.Range("D3")
$range.NumberFormat = "@"
$range.Value = Strip-MidnightTimestamp $range.Text

# --- Method 2: Cells.Item(row, column).Value -------------------------------
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = Strip-MidnightTimestamp $cell.Text

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = Strip-MidnightTimestamp $cell.Text

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = Strip-MidnightTimestamp $cell.Text
